$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.491.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "'2.899.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'566.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").Value = "'143.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").Value = "'2.897.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "'0.431"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'32.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "'3.377.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "'62.383.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'2.900.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'427.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "'13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "'6.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").Value = "'78.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'9.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").Value = "'0.0000110"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "'7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'25.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D36").Value = "'0.949"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("D37").Value = "'5.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'48.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "'2.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.05%  "
$ws.Range("D40").Value = "'1.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").Value = "'41.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").Value = "'8.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.266"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'2.720.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'133.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'0.0337"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'356.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("D50").Value = "'0.000218"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.07%  "
$ws.Range("E51").Value = "  -0.72%  "
